# Vagus Nerve model - major schema / terminology change:
#   publications -> references
#   external     -> ontologyTerms
# (as described in the commit message: "Publication -> Reference",
#  "External -> OntologyTerm")
#
# This renames the "publications" worksheet to "references" and updates
# every cell whose text contains the old terms "external" / "publications"
# so that it uses the new terms, preserving surrounding text.

$wb = $excel.ActiveWorkbook

# 1. Rename the "publications" sheet to "references"
$pubSheet = $wb.Worksheets.Item("publications")
$pubSheet.Name = "references"

# 2. Helper: replace text in a cell if it contains one of the old terms.
#    NOTE: reading `$cell.Value` (no parens) returns the property
#    descriptor in this host, so the getter must be invoked as `.Value()`.
function Update-Cell($ws, $addr, $find, $replace) {
    $cell = $ws.Range($addr)
    $old = $cell.Value()
    $cell.Value = $old.Replace($find, $replace)
}

# lyphs sheet
$ws = $wb.Worksheets.Item("lyphs")
Update-Cell $ws "C1"  "external" "ontologyTerms"
Update-Cell $ws "B16" "external" "ontologyTerms"
Update-Cell $ws "B18" "external" "ontologyTerms"

# materials sheet
$ws = $wb.Worksheets.Item("materials")
Update-Cell $ws "C1" "external" "ontologyTerms"

# links sheet
$ws = $wb.Worksheets.Item("links")
Update-Cell $ws "G1" "publications" "references"

# chains sheet
$ws = $wb.Worksheets.Item("chains")
Update-Cell $ws "C1"  "external" "ontologyTerms"
Update-Cell $ws "M1"  "publications" "references"
Update-Cell $ws "B10" "external" "ontologyTerms"
Update-Cell $ws "B11" "external" "ontologyTerms"

# nodes sheet
$ws = $wb.Worksheets.Item("nodes")
Update-Cell $ws "C1" "external" "ontologyTerms"

# neurons sheet
$ws = $wb.Worksheets.Item("neurons")
Update-Cell $ws "C6" "external" "ontologyTerms"
Update-Cell $ws "C7" "external" "ontologyTerms"
Update-Cell $ws "C8" "external" "ontologyTerms"
Update-Cell $ws "C9" "external" "ontologyTerms"

# groups sheet
$ws = $wb.Worksheets.Item("groups")
Update-Cell $ws "G1" "publications" "references"

Write-Output "done"
